# Update "想去人数" (wanted-to-go count) figures in column F across sheets,
# reflecting the regenerated data as of commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 298
$ws1.Range("F4").Value = 7966
$ws1.Range("F5").Value = 5818
$ws1.Range("F6").Value = 493
$ws1.Range("F9").Value = 69
$ws1.Range("F10").Value = 278
$ws1.Range("F11").Value = 350

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 88

# --- Sheet "全部类型" (All types, combined) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 298
$ws4.Range("F4").Value = 7966
$ws4.Range("F5").Value = 5818
$ws4.Range("F6").Value = 493
$ws4.Range("F9").Value = 69
$ws4.Range("F10").Value = 278
$ws4.Range("F11").Value = 88
$ws4.Range("F14").Value = 350
